# Auto commit at 2025-11-17 08:09:24.30
# Appends the next day's (2025-11-16, serial 45977) readings for the two
# charging stations as two new rows (154 and 155) at the bottom of the
# existing data table, then moves the active selection to D159 (matching
# where the next empty-row entry would begin two rows below the new data).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 154: 四方坪站充电量(kw) ----
$ws.Range("A154").Value = 45977
$ws.Range("B154").Value = "四方坪站充电量(kw)"
$vals154 = @(741.34500000000003, 1114.3920000000001, 562.45699999999999, 378.298, 191.78300000000002, 647.02800000000013, 450.04500000000007, 152.26399999999998, 119.02, 118.67200000000001, 126.1, 248.54800000000003, 704.72300000000018, 1105.2009999999998, 582.89900000000011, 299.07499999999999, 323.84800000000001, 350.50399999999991, 139.41200000000001, 175.71899999999999, 77.34, 72.960000000000008, 47.06, 21.8)
for ($i = 0; $i -lt $vals154.Length; $i++) {
    $ws.Cells.Item(154, 3 + $i).Value = $vals154[$i]
}

# ---- Row 155: 高岭站充电量(kw) ----
$ws.Range("A155").Value = 45977
$ws.Range("B155").Value = "高岭站充电量(kw)"
$vals155 = @(480.63400000000001, 277.923, 155.92000000000002, 75.436999999999998, 117.49299999999999, 153.39100000000002, 319.58299999999997, 42.933999999999997, 414.70800000000008, 298.20800000000008, 386.30300000000005, 76.567000000000007, 264.53100000000001, 332.57400000000007, 243.56200000000001, 307.75799999999992, 153.92700000000002, 141.10599999999999, 84.961000000000013, 23.696000000000002, 6.9109999999999996, 163.33799999999999, 15.667999999999999, 4.1059999999999999)
for ($i = 0; $i -lt $vals155.Length; $i++) {
    $ws.Cells.Item(155, 3 + $i).Value = $vals155[$i]
}

# Move the selection the same way the original author's Excel session
# ended up (two blank rows below the freshly entered data, column D).
$null = $ws.Range("D159").Select()
